$d = $word.ActiveDocument

# --- Locate the last paragraph of real content we want to keep: the one
# that ends with "...showcasing communication and public speaking skills."
$anchorText = "showcasing communication and public speaking skills."
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*$anchorText*") {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph ending in '$anchorText'"
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$lastCount  = $d.Paragraphs.Count
$lastPara   = $d.Paragraphs.Item($lastCount)

# --- Remove everything between the end of the anchor paragraph and the
# start of the very last paragraph of the document (this wipes out the
# stray blank paragraphs, the "S;dfk;..." junk paragraph, and all the
# trailing empty/indented paragraphs), while leaving the document's final
# paragraph mark (tied to the sectPr) intact so the body still ends on a
# paragraph.
if ($anchorIndex -lt $lastCount -and $lastPara.Range.Start -gt $anchorPara.Range.End) {
    $junk = $d.Range($anchorPara.Range.End, $lastPara.Range.Start)
    $junk.Delete()

    # --- The surviving final paragraph still carries the old leftover
    # direct formatting (bold/subscript/indentation/fonts) from the
    # deleted paragraphs. Clearing its style back to Normal strips all of
    # that direct paragraph/run formatting, leaving a clean, empty
    # paragraph right before the section properties.
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Style = "Normal"
}
